# Apply crypto price/volume updates per commit diff (Thu Jul 11 17:56:00 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D that would otherwise be auto-coerced to numbers by COM's .Value
# setter need an explicit Text number format first so they stay strings,
# matching the inlineStr cell type in the source workbook.
$textForceCells = @(
    "D5", "D6", "D10", "D15", "D16", "D19", "D22", "D23", "D24", "D25", "D27", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D41", "D43", "D45", "D46", "D49", "D50", "D51"
)
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "57.884.63"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "3.134.87"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "528.92"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").Value = "138.32"
$ws.Range("E6").Value = "  -2.02%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.133.48"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("E9").Value = "  +2.97%  "
$ws.Range("D10").Value = "7.20"
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("E12").Value = "  +2.55%  "
$ws.Range("D13").Value = "3.674.75"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("E14").Value = "  +2.37%  "
$ws.Range("D15").Value = "25.51"
$ws.Range("E15").Value = "  -2.83%  "
$ws.Range("D16").Value = "0.0000165"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("D17").Value = "58.013.13"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").Value = "3.134.85"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").Value = "6.00"
$ws.Range("E19").Value = "  -2.14%  "
$ws.Range("E20").Value = "  -0.55%  "
$ws.Range("E21").Value = "  -1.59%  "
$ws.Range("D22").Value = "354.29"
$ws.Range("E22").Value = "  +5.28%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "68.94"
$ws.Range("E24").Value = "  +3.57%  "
$ws.Range("D25").Value = "0.508"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").Value = "0.0₃0919"
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("D29").Value = "7.52"
$ws.Range("E29").Value = "  +4.27%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "6.19"
$ws.Range("E31").Value = "  -5.52%  "
$ws.Range("D32").Value = "1.89"
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("D33").Value = "21.21"
$ws.Range("E33").Value = "  +0.89%  "
$ws.Range("D34").Value = "1.19"
$ws.Range("E34").Value = "  -1.23%  "
$ws.Range("D35").Value = "4.99"
$ws.Range("E35").Value = "  +6.87%  "
$ws.Range("D36").Value = "158.78"
$ws.Range("E36").Value = "  +2.42%  "
$ws.Range("E37").Value = "  +0.91%  "
$ws.Range("D38").Value = "26.70"
$ws.Range("E38").Value = "  -1.39%  "
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("D41").Value = "4.19"
$ws.Range("E41").Value = "  +6.50%  "
$ws.Range("E42").Value = "  +5.98%  "
$ws.Range("D43").Value = "0.703"
$ws.Range("E43").Value = "  +2.37%  "
$ws.Range("D44").Value = "3.175.46"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("D45").Value = "0.0272"
$ws.Range("E45").Value = "  +4.88%  "
$ws.Range("D46").Value = "36.60"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").Value = "2.306.82"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("D49").Value = "0.970"
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("D50").Value = "6.04"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").Value = "20.38"
$ws.Range("E51").Value = "  -1.86%  "

Write-Output "Applied cryptos list update"
